$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force BF column (Date, stored as text like "2012-04-01") to remain text,
# so Excel does not auto-convert the ISO-looking string into a date serial number.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Row 2
$ws.Range("AG2").Value = 9
$ws.Range("AI2").Value = 17
$ws.Range("AK2").Value = 16
$ws.Range("AR2").Value = 25
$ws.Range("AS2").Value = 11
$ws.Range("AU2").Value = 9
$ws.Range("AW2").Value = 12
$ws.Range("AX2").Value = 23
$ws.Range("BC2").Value = 11
$ws.Range("BF2").Value = "2012-04-01"

# Row 3
$ws.Range("D3").Value = 51
$ws.Range("E3").Value = 29
$ws.Range("G3").Value = 0.569
$ws.Range("I3").Value = 35.3
$ws.Range("J3").Value = 77.09999999999999
$ws.Range("K3").Value = 0.457
$ws.Range("M3").Value = 15.1
$ws.Range("N3").Value = 0.359
$ws.Range("O3").Value = 15.6
$ws.Range("Q3").Value = 0.778
$ws.Range("S3").Value = 30.4
$ws.Range("T3").Value = 38.4
$ws.Range("Y3").Value = 4.6
$ws.Range("AA3").Value = 18.9
$ws.Range("AC3").Value = 1.4
$ws.Range("AD3").Value = 20
$ws.Range("AE3").Value = 11
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 19
$ws.Range("AL3").Value = 23
$ws.Range("AM3").Value = 23
$ws.Range("AO3").Value = 21
$ws.Range("AP3").Value = 27
$ws.Range("AQ3").Value = 6
$ws.Range("AV3").Value = 16
$ws.Range("AX3").Value = 5
$ws.Range("AY3").Value = 9
$ws.Range("BC3").Value = 13
$ws.Range("BF3").Value = "2012-04-01"

# Row 4
$ws.Range("AD4").Value = 25
$ws.Range("AQ4").Value = 16
$ws.Range("AR4").Value = 26
$ws.Range("AX4").Value = 6
$ws.Range("AY4").Value = 25
$ws.Range("BA4").Value = 14
$ws.Range("BF4").Value = "2012-04-01"

# Row 5
$ws.Range("D5").Value = 53
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 0.792
$ws.Range("I5").Value = 37.7
$ws.Range("J5").Value = 82.40000000000001
$ws.Range("K5").Value = 0.457
$ws.Range("M5").Value = 16.5
$ws.Range("N5").Value = 0.381
$ws.Range("O5").Value = 15.6
$ws.Range("P5").Value = 21.5
$ws.Range("Q5").Value = 0.727
$ws.Range("R5").Value = 13.7
$ws.Range("S5").Value = 32.3
$ws.Range("T5").Value = 46
$ws.Range("U5").Value = 23.1
$ws.Range("W5").Value = 7
$ws.Range("X5").Value = 5.9
$ws.Range("Y5").Value = 5.1
$ws.Range("Z5").Value = 17.2
$ws.Range("AA5").Value = 17.7
$ws.Range("AB5").Value = 97.2
$ws.Range("AC5").Value = 8.699999999999999
$ws.Range("AD5").Value = 3
$ws.Range("AI5").Value = 6
$ws.Range("AJ5").Value = 10
$ws.Range("AK5").Value = 6
$ws.Range("AL5").Value = 18
$ws.Range("AO5").Value = 22
$ws.Range("AP5").Value = 20
$ws.Range("AQ5").Value = 25
$ws.Range("AV5").Value = 5
$ws.Range("AY5").Value = 20
$ws.Range("BB5").Value = 11
$ws.Range("BF5").Value = "2012-04-01"

# Row 6
$ws.Range("AD6").Value = 25
$ws.Range("AE6").Value = 26
$ws.Range("AH6").Value = 17
$ws.Range("AP6").Value = 5
$ws.Range("BF6").Value = "2012-04-01"

# Row 7
$ws.Range("AD7").Value = 3
$ws.Range("AE7").Value = 9
$ws.Range("AH7").Value = 15
$ws.Range("AI7").Value = 19
$ws.Range("AJ7").Value = 14
$ws.Range("AO7").Value = 23
$ws.Range("AX7").Value = 14
$ws.Range("BF7").Value = "2012-04-01"

# Row 8
$ws.Range("D8").Value = 52
$ws.Range("E8").Value = 28
$ws.Range("G8").Value = 0.538
$ws.Range("J8").Value = 81.40000000000001
$ws.Range("K8").Value = 0.47
$ws.Range("N8").Value = 0.321
$ws.Range("O8").Value = 20.5
$ws.Range("P8").Value = 27.7
$ws.Range("R8").Value = 11.1
$ws.Range("T8").Value = 43.6
$ws.Range("U8").Value = 23.6
$ws.Range("W8").Value = 8.300000000000001
$ws.Range("X8").Value = 5.2
$ws.Range("AA8").Value = 22.6
$ws.Range("AC8").Value = 1.4
$ws.Range("AD8").Value = 8
$ws.Range("AE8").Value = 13
$ws.Range("AL8").Value = 16
$ws.Range("AR8").Value = 17
$ws.Range("AW8").Value = 6
$ws.Range("AZ8").Value = 14
$ws.Range("BA8").Value = 1
$ws.Range("BB8").Value = 2
$ws.Range("BC8").Value = 12
$ws.Range("BF8").Value = "2012-04-01"

# Row 9
$ws.Range("AD9").Value = 8
$ws.Range("AH9").Value = 12
$ws.Range("AP9").Value = 18
$ws.Range("BF9").Value = "2012-04-01"

# Row 10
$ws.Range("D10").Value = 50
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 0.4
$ws.Range("I10").Value = 37.2
$ws.Range("J10").Value = 81.59999999999999
$ws.Range("K10").Value = 0.456
$ws.Range("N10").Value = 0.384
$ws.Range("O10").Value = 15.3
$ws.Range("P10").Value = 19.7
$ws.Range("Q10").Value = 0.775
$ws.Range("S10").Value = 29.7
$ws.Range("T10").Value = 39.3
$ws.Range("U10").Value = 22.1
$ws.Range("V10").Value = 13.9
$ws.Range("Y10").Value = 4.5
$ws.Range("AB10").Value = 97.8
$ws.Range("AC10").Value = -2.2
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 11
$ws.Range("AJ10").Value = 13
$ws.Range("AK10").Value = 7
$ws.Range("AL10").Value = 4
$ws.Range("AO10").Value = 26
$ws.Range("AQ10").Value = 8
$ws.Range("AS10").Value = 23
$ws.Range("AV10").Value = 6
$ws.Range("AY10").Value = 7
$ws.Range("BF10").Value = "2012-04-01"

# Row 11
$ws.Range("D11").Value = 52
$ws.Range("F11").Value = 24
$ws.Range("G11").Value = 0.538
$ws.Range("H11").Value = 48.7
$ws.Range("I11").Value = 37.6
$ws.Range("J11").Value = 83
$ws.Range("L11").Value = 6.9
$ws.Range("M11").Value = 19.4
$ws.Range("N11").Value = 0.357
$ws.Range("O11").Value = 15.9
$ws.Range("P11").Value = 20.2
$ws.Range("Q11").Value = 0.788
$ws.Range("R11").Value = 11.4
$ws.Range("S11").Value = 30.8
$ws.Range("U11").Value = 20.9
$ws.Range("V11").Value = 14.9
$ws.Range("Y11").Value = 5.1
$ws.Range("AA11").Value = 18.5
$ws.Range("AD11").Value = 8
$ws.Range("AE11").Value = 13
$ws.Range("AF11").Value = 14
$ws.Range("AG11").Value = 14
$ws.Range("AH11").Value = 5
$ws.Range("AI11").Value = 7
$ws.Range("AJ11").Value = 6
$ws.Range("AO11").Value = 19
$ws.Range("AP11").Value = 26
$ws.Range("AQ11").Value = 3
$ws.Range("AU11").Value = 16
$ws.Range("AV11").Value = 17
$ws.Range("AX11").Value = 18
$ws.Range("AY11").Value = 19
$ws.Range("AZ11").Value = 20
$ws.Range("BF11").Value = "2012-04-01"

# Row 12
$ws.Range("D12").Value = 51
$ws.Range("E12").Value = 30
$ws.Range("G12").Value = 0.588
$ws.Range("H12").Value = 48.3
$ws.Range("L12").Value = 5.7
$ws.Range("N12").Value = 0.369
$ws.Range("O12").Value = 20.1
$ws.Range("Q12").Value = 0.776
$ws.Range("R12").Value = 12.1
$ws.Range("S12").Value = 31
$ws.Range("T12").Value = 43.1
$ws.Range("U12").Value = 18.2
$ws.Range("V12").Value = 14.3
$ws.Range("W12").Value = 8
$ws.Range("X12").Value = 5.3
$ws.Range("Y12").Value = 5.9
$ws.Range("Z12").Value = 22
$ws.Range("AB12").Value = 96.2
$ws.Range("AD12").Value = 20
$ws.Range("AE12").Value = 9
$ws.Range("AF12").Value = 7
$ws.Range("AG12").Value = 8
$ws.Range("AH12").Value = 19
$ws.Range("AK12").Value = 23
$ws.Range("AQ12").Value = 7
$ws.Range("AR12").Value = 9
$ws.Range("AS12").Value = 12
$ws.Range("AX12").Value = 11
$ws.Range("BF12").Value = "2012-04-01"

# Row 13
$ws.Range("AD13").Value = 8
$ws.Range("AF13").Value = 7
$ws.Range("AS13").Value = 26
$ws.Range("AW13").Value = 11
$ws.Range("AX13").Value = 22
$ws.Range("BB13").Value = 12
$ws.Range("BF13").Value = "2012-04-01"

# Row 14
$ws.Range("D14").Value = 52
$ws.Range("E14").Value = 32
$ws.Range("G14").Value = 0.615
$ws.Range("I14").Value = 36.1
$ws.Range("J14").Value = 79.40000000000001
$ws.Range("K14").Value = 0.454
$ws.Range("L14").Value = 5.4
$ws.Range("N14").Value = 0.313
$ws.Range("O14").Value = 18.3
$ws.Range("P14").Value = 24.2
$ws.Range("Q14").Value = 0.758
$ws.Range("T14").Value = 45.6
$ws.Range("U14").Value = 21.7
$ws.Range("V14").Value = 15.3
$ws.Range("AA14").Value = 20.4
$ws.Range("AB14").Value = 95.8
$ws.Range("AC14").Value = 2.3
$ws.Range("AD14").Value = 8
$ws.Range("AI14").Value = 18
$ws.Range("AK14").Value = 9
$ws.Range("AL14").Value = 24
$ws.Range("AN14").Value = 28
$ws.Range("AQ14").Value = 14
$ws.Range("AR14").Value = 14
$ws.Range("AU14").Value = 10
$ws.Range("BF14").Value = "2012-04-01"

# Row 15
$ws.Range("AD15").Value = 25
$ws.Range("AE15").Value = 13
$ws.Range("AK15").Value = 15
$ws.Range("AN15").Value = 25
$ws.Range("AR15").Value = 7
$ws.Range("BF15").Value = "2012-04-01"

# Row 16
$ws.Range("D16").Value = 50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 0.74
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 78.90000000000001
$ws.Range("K16").Value = 0.482
$ws.Range("L16").Value = 5.7
$ws.Range("M16").Value = 15
$ws.Range("N16").Value = 0.378
$ws.Range("O16").Value = 19.6
$ws.Range("P16").Value = 25.2
$ws.Range("Q16").Value = 0.779
$ws.Range("R16").Value = 10.1
$ws.Range("T16").Value = 41.7
$ws.Range("U16").Value = 20.9
$ws.Range("V16").Value = 15
$ws.Range("W16").Value = 9.1
$ws.Range("Z16").Value = 20
$ws.Range("AA16").Value = 20.9
$ws.Range("AB16").Value = 101.3
$ws.Range("AC16").Value = 7.7
$ws.Range("AH16").Value = 4
$ws.Range("AJ16").Value = 26
$ws.Range("AM16").Value = 24
$ws.Range("AQ16").Value = 5
$ws.Range("AR16").Value = 27
$ws.Range("AT16").Value = 19
$ws.Range("AU16").Value = 15
$ws.Range("AX16").Value = 10
$ws.Range("AY16").Value = 4
$ws.Range("BA16").Value = 7
$ws.Range("BF16").Value = "2012-04-01"

# Row 17
$ws.Range("AD17").Value = 8
$ws.Range("AE17").Value = 20
$ws.Range("AG17").Value = 20
$ws.Range("AH17").Value = 27
$ws.Range("AI17").Value = 8
$ws.Range("AQ17").Value = 4
$ws.Range("AR17").Value = 6
$ws.Range("AT17").Value = 18
$ws.Range("AV17").Value = 4
$ws.Range("AW17").Value = 9
$ws.Range("AY17").Value = 10
$ws.Range("BF17").Value = "2012-04-01"

# Row 18
$ws.Range("D18").Value = 53
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 0.472
$ws.Range("L18").Value = 7.2
$ws.Range("N18").Value = 0.337
$ws.Range("O18").Value = 19.6
$ws.Range("P18").Value = 25.4
$ws.Range("Q18").Value = 0.773
$ws.Range("R18").Value = 12.6
$ws.Range("S18").Value = 32
$ws.Range("T18").Value = 44.6
$ws.Range("U18").Value = 19.5
$ws.Range("V18").Value = 15.5
$ws.Range("AA18").Value = 22
$ws.Range("AB18").Value = 98.59999999999999
$ws.Range("AC18").Value = -0.2
$ws.Range("AD18").Value = 3
$ws.Range("AE18").Value = 18
$ws.Range("AF18").Value = 19
$ws.Range("AG18").Value = 19
$ws.Range("AH18").Value = 15
$ws.Range("AJ18").Value = 9
$ws.Range("AP18").Value = 6
$ws.Range("AQ18").Value = 9
$ws.Range("AR18").Value = 5
$ws.Range("AV18").Value = 25
$ws.Range("BF18").Value = "2012-04-01"

# Row 19
$ws.Range("AH19").Value = 29
$ws.Range("AL19").Value = 3
$ws.Range("AR19").Value = 8
$ws.Range("BF19").Value = "2012-04-01"

# Row 20
$ws.Range("D20").Value = 52
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 0.25
$ws.Range("I20").Value = 35.1
$ws.Range("M20").Value = 11.8
$ws.Range("N20").Value = 0.33
$ws.Range("O20").Value = 15.1
$ws.Range("P20").Value = 19.9
$ws.Range("R20").Value = 11.3
$ws.Range("V20").Value = 15.4
$ws.Range("W20").Value = 7.1
$ws.Range("X20").Value = 4.7
$ws.Range("Y20").Value = 5.9
$ws.Range("AA20").Value = 18.3
$ws.Range("AB20").Value = 89.2
$ws.Range("AC20").Value = -4.5
$ws.Range("AD20").Value = 8
$ws.Range("AH20").Value = 21
$ws.Range("AK20").Value = 14
$ws.Range("AQ20").Value = 15
$ws.Range("AS20").Value = 21
$ws.Range("AV20").Value = 23
$ws.Range("AW20").Value = 22
$ws.Range("AX20").Value = 20
$ws.Range("AY20").Value = 26
$ws.Range("BF20").Value = "2012-04-01"

# Row 21
$ws.Range("AD21").Value = 3
$ws.Range("AH21").Value = 23
$ws.Range("AK21").Value = 22
$ws.Range("AN21").Value = 26
$ws.Range("AR21").Value = 13
$ws.Range("BA21").Value = 2
$ws.Range("BF21").Value = "2012-04-01"

# Row 22
$ws.Range("D22").Value = 51
$ws.Range("E22").Value = 39
$ws.Range("G22").Value = 0.765
$ws.Range("J22").Value = 78.90000000000001
$ws.Range("M22").Value = 20.1
$ws.Range("N22").Value = 0.358
$ws.Range("O22").Value = 21
$ws.Range("P22").Value = 26.4
$ws.Range("Q22").Value = 0.797
$ws.Range("R22").Value = 10.9
$ws.Range("S22").Value = 32.6
$ws.Range("T22").Value = 43.5
$ws.Range("U22").Value = 18.5
$ws.Range("V22").Value = 16.4
$ws.Range("Z22").Value = 20.3
$ws.Range("AA22").Value = 20.1
$ws.Range("AB22").Value = 103.7
$ws.Range("AC22").Value = 6.6
$ws.Range("AD22").Value = 20
$ws.Range("AF22").Value = 2
$ws.Range("AI22").Value = 5
$ws.Range("AJ22").Value = 25
$ws.Range("AM22").Value = 11
$ws.Range("AT22").Value = 6
$ws.Range("AY22").Value = 11
$ws.Range("AZ22").Value = 22
$ws.Range("BA22").Value = 13
$ws.Range("BB22").Value = 1
$ws.Range("BF22").Value = "2012-04-01"

# Row 23
$ws.Range("D23").Value = 52
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 0.615
$ws.Range("I23").Value = 34.3
$ws.Range("J23").Value = 77.3
$ws.Range("K23").Value = 0.444
$ws.Range("L23").Value = 10.3
$ws.Range("M23").Value = 26.9
$ws.Range("N23").Value = 0.383
$ws.Range("P23").Value = 23.7
$ws.Range("Q23").Value = 0.643
$ws.Range("R23").Value = 10.9
$ws.Range("S23").Value = 31.9
$ws.Range("U23").Value = 20.3
$ws.Range("V23").Value = 15.6
$ws.Range("X23").Value = 4.3
$ws.Range("Z23").Value = 17.8
$ws.Range("AB23").Value = 94.2
$ws.Range("AC23").Value = 1.8
$ws.Range("AD23").Value = 8
$ws.Range("AE23").Value = 5
$ws.Range("AF23").Value = 5
$ws.Range("AG23").Value = 5
$ws.Range("AH23").Value = 12
$ws.Range("AO23").Value = 27
$ws.Range("AR23").Value = 18
$ws.Range("BC23").Value = 10
$ws.Range("BF23").Value = "2012-04-01"

# Row 24
$ws.Range("AD24").Value = 8
$ws.Range("AE24").Value = 11
$ws.Range("AI24").Value = 8
$ws.Range("AT24").Value = 8
$ws.Range("AU24").Value = 8
$ws.Range("AX24").Value = 17
$ws.Range("BF24").Value = "2012-04-01"

# Row 25
$ws.Range("D25").Value = 51
$ws.Range("E25").Value = 25
$ws.Range("G25").Value = 0.49
$ws.Range("I25").Value = 37.2
$ws.Range("J25").Value = 81.7
$ws.Range("K25").Value = 0.455
$ws.Range("M25").Value = 19.2
$ws.Range("N25").Value = 0.342
$ws.Range("O25").Value = 15.8
$ws.Range("Q25").Value = 0.753
$ws.Range("S25").Value = 30.8
$ws.Range("V25").Value = 14.1
$ws.Range("X25").Value = 5.5
$ws.Range("AC25").Value = -0.9
$ws.Range("AD25").Value = 20
$ws.Range("AJ25").Value = 11
$ws.Range("AL25").Value = 15
$ws.Range("AN25").Value = 17
$ws.Range("AO25").Value = 20
$ws.Range("AQ25").Value = 17
$ws.Range("AS25").Value = 15
$ws.Range("AU25").Value = 5
$ws.Range("AV25").Value = 9
$ws.Range("AX25").Value = 9
$ws.Range("AY25").Value = 5
$ws.Range("BF25").Value = "2012-04-01"

# Row 26
$ws.Range("D26").Value = 52
$ws.Range("E26").Value = 24
$ws.Range("G26").Value = 0.462
$ws.Range("I26").Value = 36.3
$ws.Range("J26").Value = 81.59999999999999
$ws.Range("K26").Value = 0.444
$ws.Range("L26").Value = 6.9
$ws.Range("M26").Value = 20.3
$ws.Range("Q26").Value = 0.79
$ws.Range("R26").Value = 10.9
$ws.Range("T26").Value = 40.6
$ws.Range("V26").Value = 14.2
$ws.Range("W26").Value = 8.199999999999999
$ws.Range("AA26").Value = 20.3
$ws.Range("AB26").Value = 96.8
$ws.Range("AC26").Value = 0.8
$ws.Range("AD26").Value = 8
$ws.Range("AE26").Value = 20
$ws.Range("AG26").Value = 20
$ws.Range("AH26").Value = 12
$ws.Range("AK26").Value = 17
$ws.Range("AN26").Value = 18
$ws.Range("AP26").Value = 17
$ws.Range("AR26").Value = 18
$ws.Range("AS26").Value = 25
$ws.Range("AV26").Value = 11
$ws.Range("AW26").Value = 10
$ws.Range("BB26").Value = 14
$ws.Range("BF26").Value = "2012-04-01"

# Row 27
$ws.Range("AD27").Value = 8
$ws.Range("AH27").Value = 21
$ws.Range("AI27").Value = 10
$ws.Range("AL27").Value = 17
$ws.Range("AM27").Value = 12
$ws.Range("AN27").Value = 27
$ws.Range("AT27").Value = 7
$ws.Range("BA27").Value = 9
$ws.Range("BF27").Value = "2012-04-01"

# Row 28
$ws.Range("AD28").Value = 25
$ws.Range("AF28").Value = 4
$ws.Range("AJ28").Value = 7
$ws.Range("AO28").Value = 24
$ws.Range("AP28").Value = 21
$ws.Range("AQ28").Value = 26
$ws.Range("AU28").Value = 6
$ws.Range("AW28").Value = 21
$ws.Range("BF28").Value = "2012-04-01"

# Row 29
$ws.Range("D29").Value = 52
$ws.Range("E29").Value = 17
$ws.Range("G29").Value = 0.327
$ws.Range("I29").Value = 34.6
$ws.Range("J29").Value = 78.5
$ws.Range("K29").Value = 0.441
$ws.Range("N29").Value = 0.336
$ws.Range("O29").Value = 16.9
$ws.Range("Q29").Value = 0.771
$ws.Range("T29").Value = 41.3
$ws.Range("U29").Value = 21.3
$ws.Range("V29").Value = 15.3
$ws.Range("Z29").Value = 23.8
$ws.Range("AA29").Value = 18.7
$ws.Range("AB29").Value = 91.7
$ws.Range("AC29").Value = -3.9
$ws.Range("AD29").Value = 8
$ws.Range("AE29").Value = 26
$ws.Range("AH29").Value = 7
$ws.Range("AL29").Value = 21
$ws.Range("AQ29").Value = 10
$ws.Range("AS29").Value = 14
$ws.Range("AV29").Value = 20
$ws.Range("BF29").Value = "2012-04-01"

# Row 30
$ws.Range("AD30").Value = 3
$ws.Range("AS30").Value = 17
$ws.Range("AW30").Value = 8
$ws.Range("AY30").Value = 24
$ws.Range("BA30").Value = 8
$ws.Range("BF30").Value = "2012-04-01"

# Row 31
$ws.Range("D31").Value = 51
$ws.Range("F31").Value = 39
$ws.Range("G31").Value = 0.235
$ws.Range("J31").Value = 83.09999999999999
$ws.Range("K31").Value = 0.436
$ws.Range("L31").Value = 5.1
$ws.Range("N31").Value = 0.317
$ws.Range("Q31").Value = 0.721
$ws.Range("R31").Value = 11.9
$ws.Range("S31").Value = 30
$ws.Range("T31").Value = 41.9
$ws.Range("U31").Value = 18.3
$ws.Range("W31").Value = 7.9
$ws.Range("AA31").Value = 18.7
$ws.Range("AD31").Value = 20
$ws.Range("AI31").Value = 16
$ws.Range("AJ31").Value = 5
$ws.Range("AO31").Value = 25
$ws.Range("AS31").Value = 20
$ws.Range("AT31").Value = 17
$ws.Range("AY31").Value = 8
$ws.Range("BF31").Value = "2012-04-01"
